$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.7667715243636337
$ws.Cells.Item(2, 3).Value = 0.2713778478191671
$ws.Cells.Item(2, 4).Value = 0.08107532069008982
$ws.Cells.Item(2, 5).Value = 0.1350018207781183
$ws.Cells.Item(2, 6).Value = 2.72955882362993
$ws.Cells.Item(2, 9).Value = 1.369635133289705
$ws.Cells.Item(2, 10).Value = 0.2156021860236734
$ws.Cells.Item(2, 11).Value = 1.030382707340522
$ws.Cells.Item(2, 14).Value = 2.548320960721064
$ws.Cells.Item(3, 2).Value = 0.7263890224822944
$ws.Cells.Item(3, 3).Value = 0.2603134060771026
$ws.Cells.Item(3, 4).Value = 0.07855395878269888
$ws.Cells.Item(3, 5).Value = 0.1318959325522293
$ws.Cells.Item(3, 6).Value = 2.717530195679089
$ws.Cells.Item(3, 9).Value = 1.3692536558898
$ws.Cells.Item(3, 10).Value = 0.2113600284143615
$ws.Cells.Item(3, 11).Value = 0.9797632597496602
$ws.Cells.Item(3, 14).Value = 2.564830886544549
$ws.Cells.Item(4, 2).Value = 0.702025315822624
$ws.Cells.Item(4, 3).Value = 0.253684088700453
$ws.Cells.Item(4, 4).Value = 0.07704666928807313
$ws.Cells.Item(4, 5).Value = 0.1300628855400099
$ws.Cells.Item(4, 6).Value = 2.711447463747049
$ws.Cells.Item(4, 9).Value = 1.369619534222224
$ws.Cells.Item(4, 10).Value = 0.2088774292418734
$ws.Cells.Item(4, 11).Value = 0.9492766824182581
$ws.Cells.Item(4, 14).Value = 2.575691282604836
$ws.Cells.Item(5, 2).Value = 0.6922053995449744
$ws.Cells.Item(5, 3).Value = 0.2510237924616376
$ws.Cells.Item(5, 4).Value = 0.07644274143548557
$ws.Cells.Item(5, 5).Value = 0.1293345034873283
$ws.Cells.Item(5, 6).Value = 2.709296048314187
$ws.Cells.Item(5, 9).Value = 1.369919395441833
$ws.Cells.Item(5, 10).Value = 0.207896404893674
$ws.Cells.Item(5, 11).Value = 0.9370024204892502
$ws.Cells.Item(5, 14).Value = 2.580298677090624
$ws.Cells.Item(6, 2).Value = 0.6905813637756921
$ws.Cells.Item(6, 3).Value = 0.2505845378130402
$ws.Cells.Item(6, 4).Value = 0.07634308313484439
$ws.Cells.Item(6, 5).Value = 0.1292146794776521
$ws.Cells.Item(6, 6).Value = 2.708958570481656
$ws.Cells.Item(6, 9).Value = 1.369978288872012
$ws.Cells.Item(6, 10).Value = 0.2077353566308489
$ws.Cells.Item(6, 11).Value = 0.9349733014612411
$ws.Cells.Item(6, 14).Value = 2.581074700493801
$ws.Cells.Item(7, 2).Value = 0.7018924416694858
$ws.Cells.Item(7, 3).Value = 0.2536480443220057
$ws.Cells.Item(7, 4).Value = 0.07703848273277458
$ws.Cells.Item(7, 5).Value = 0.1300529870201252
$ws.Cells.Item(7, 6).Value = 2.711417123964438
$ws.Cells.Item(7, 9).Value = 1.36962296801287
$ws.Cells.Item(7, 10).Value = 0.2088640747466641
$ws.Cells.Item(7, 11).Value = 0.9491105431764311
$ws.Cells.Item(7, 14).Value = 2.575752684098724
$ws.Cells.Item(8, 2).Value = 0.7527581237681886
$ws.Cells.Item(8, 3).Value = 0.2675286589344523
$ws.Cells.Item(8, 4).Value = 0.08019749994583236
$ws.Cells.Item(8, 5).Value = 0.1339155521335051
$ws.Cells.Item(8, 6).Value = 2.72514075881476
$ws.Cells.Item(8, 9).Value = 1.369378960174068
$ws.Cells.Item(8, 10).Value = 0.2141141162422926
$ws.Cells.Item(8, 11).Value = 1.0128057680175
$ws.Cells.Item(8, 14).Value = 2.553863374647563
$ws.Cells.Item(9, 2).Value = 0.8559326802432565
$ws.Cells.Item(9, 3).Value = 0.2960588165684896
$ws.Cells.Item(9, 4).Value = 0.08671518670705325
$ws.Cells.Item(9, 5).Value = 0.1420777950149841
$ws.Cells.Item(9, 6).Value = 2.762407981765648
$ws.Cells.Item(9, 9).Value = 1.373669887941404
$ws.Cells.Item(9, 10).Value = 0.2253812985056811
$ws.Cells.Item(9, 11).Value = 1.142437021405044
$ws.Cells.Item(9, 14).Value = 2.516684461755077
$ws.Cells.Item(10, 2).Value = 0.9338402286476537
$ws.Cells.Item(10, 3).Value = 0.3178306191757656
$ws.Cells.Item(10, 4).Value = 0.09169963031925477
$ws.Cells.Item(10, 5).Value = 0.14843482373076
$ws.Cells.Item(10, 6).Value = 2.796132083912823
$ws.Cells.Item(10, 9).Value = 1.379743150831558
$ws.Cells.Item(10, 10).Value = 0.2342571856878948
$ws.Cells.Item(10, 11).Value = 1.240586262868902
$ws.Cells.Item(10, 14).Value = 2.492882011755796
$ws.Cells.Item(11, 2).Value = 0.9697435849798808
$ws.Cells.Item(11, 3).Value = 0.3279139580951664
$ws.Cells.Item(11, 4).Value = 0.0940095702992636
$ws.Cells.Item(11, 5).Value = 0.1514054942526499
$ws.Cells.Item(11, 6).Value = 2.812858823555317
$ws.Cells.Item(11, 9).Value = 1.383143276038645
$ws.Cells.Item(11, 10).Value = 0.2384261352473516
$ws.Cells.Item(11, 11).Value = 1.285875558275734
$ws.Cells.Item(11, 14).Value = 2.482818692443018
$ws.Cells.Item(12, 2).Value = 0.9834059455107536
$ws.Cells.Item(12, 3).Value = 0.3317582010879789
$ws.Cells.Item(12, 4).Value = 0.09489036937600304
$ws.Cells.Item(12, 5).Value = 0.1525417665050597
$ws.Cells.Item(12, 6).Value = 2.819392504738801
$ws.Cells.Item(12, 9).Value = 1.384522680477438
$ws.Cells.Item(12, 10).Value = 0.2400237617334398
$ws.Cells.Item(12, 11).Value = 1.303117896873715
$ws.Cells.Item(12, 14).Value = 2.479118127891908
$ws.Cells.Item(13, 2).Value = 0.9804605525689567
$ws.Cells.Item(13, 3).Value = 0.3309291205556519
$ws.Cells.Item(13, 4).Value = 0.09470040390854706
$ws.Cells.Item(13, 5).Value = 0.1522965452825744
$ws.Cells.Item(13, 6).Value = 2.817976474662629
$ws.Cells.Item(13, 9).Value = 1.384221512768036
$ws.Cells.Item(13, 10).Value = 0.2396788411431174
$ws.Cells.Item(13, 11).Value = 1.29940034505546
$ws.Cells.Item(13, 14).Value = 2.479910205533642
$ws.Cells.Item(14, 2).Value = 0.9708662621339386
$ws.Cells.Item(14, 3).Value = 0.3282297064395152
$ws.Cells.Item(14, 4).Value = 0.09408191261897514
$ws.Cells.Item(14, 5).Value = 0.1514987486563371
$ws.Cells.Item(14, 6).Value = 2.813392350290627
$ws.Cells.Item(14, 9).Value = 1.383254918498388
$ws.Cells.Item(14, 10).Value = 0.2385571930911539
$ws.Cells.Item(14, 11).Value = 1.287292245936385
$ws.Cells.Item(14, 14).Value = 2.482512034789266
$ws.Cells.Item(15, 2).Value = 0.9649981471709452
$ws.Cells.Item(15, 3).Value = 0.3265796142908641
$ws.Cells.Item(15, 4).Value = 0.09370385894199273
$ws.Cells.Item(15, 5).Value = 0.1510115527514202
$ws.Cells.Item(15, 6).Value = 2.810610451140988
$ws.Cells.Item(15, 9).Value = 1.382674819222686
$ws.Cells.Item(15, 10).Value = 0.2378726190697193
$ws.Cells.Item(15, 11).Value = 1.279887705305839
$ws.Cells.Item(15, 14).Value = 2.484120088820916
$ws.Cells.Item(16, 2).Value = 0.9315031716763258
$ws.Cells.Item(16, 3).Value = 0.3171752668037016
$ws.Cells.Item(16, 4).Value = 0.09154952232157143
$ws.Cells.Item(16, 5).Value = 0.1482422703407451
$ws.Cells.Item(16, 6).Value = 2.795066865316301
$ws.Cells.Item(16, 9).Value = 1.379533786977689
$ws.Cells.Item(16, 10).Value = 0.2339873815182045
$ws.Cells.Item(16, 11).Value = 1.237639404777127
$ws.Cells.Item(16, 14).Value = 2.493555079581981
$ws.Cells.Item(17, 2).Value = 0.9110736169455436
$ws.Cells.Item(17, 3).Value = 0.3114520009980879
$ws.Cells.Item(17, 4).Value = 0.0902387644421907
$ws.Cells.Item(17, 5).Value = 0.1465636005731739
$ws.Cells.Item(17, 6).Value = 2.785886513406794
$ws.Cells.Item(17, 9).Value = 1.377770250478903
$ws.Cells.Item(17, 10).Value = 0.2316375764883389
$ws.Cells.Item(17, 11).Value = 1.211885634241298
$ws.Cells.Item(17, 14).Value = 2.499539137773127
$ws.Cells.Item(18, 2).Value = 0.8993666409330103
$ws.Cells.Item(18, 3).Value = 0.3081769927461551
$ws.Cells.Item(18, 4).Value = 0.08948885277121121
$ws.Cells.Item(18, 5).Value = 0.1456054930734965
$ws.Cells.Item(18, 6).Value = 2.780736601070444
$ws.Cells.Item(18, 9).Value = 1.376815890890114
$ws.Cells.Item(18, 10).Value = 0.2302983823195319
$ws.Cells.Item(18, 11).Value = 1.197133045249331
$ws.Cells.Item(18, 14).Value = 2.50305295567064
$ws.Cells.Item(19, 2).Value = 0.8954103426145537
$ws.Cells.Item(19, 3).Value = 0.307071023591277
$ws.Cells.Item(19, 4).Value = 0.08923563392389156
$ws.Cells.Item(19, 5).Value = 0.145282368292726
$ws.Cells.Item(19, 6).Value = 2.779015310708061
$ws.Cells.Item(19, 9).Value = 1.37650305677392
$ws.Cells.Item(19, 10).Value = 0.2298470737540299
$ws.Cells.Item(19, 11).Value = 1.192148424560003
$ws.Cells.Item(19, 14).Value = 2.504255023080432
$ws.Cells.Item(20, 2).Value = 0.913243869381887
$ws.Cells.Item(20, 3).Value = 0.3120595062557641
$ws.Cells.Item(20, 4).Value = 0.09037788296384974
$ws.Cells.Item(20, 5).Value = 0.1467415299418064
$ws.Cells.Item(20, 6).Value = 2.786850280489219
$ws.Cells.Item(20, 9).Value = 1.37795177281626
$ws.Cells.Item(20, 10).Value = 0.2318864384757262
$ws.Cells.Item(20, 11).Value = 1.214620925465539
$ws.Cells.Item(20, 14).Value = 2.498894676571283
$ws.Cells.Item(21, 2).Value = 0.9736825320528624
$ws.Cells.Item(21, 3).Value = 0.3290218858193157
$ws.Cells.Item(21, 4).Value = 0.09426341391073834
$ws.Cells.Item(21, 5).Value = 0.1517327729250653
$ws.Cells.Item(21, 6).Value = 2.814733397683369
$ws.Cells.Item(21, 9).Value = 1.383536336265003
$ws.Cells.Item(21, 10).Value = 0.2388861339176742
$ws.Cells.Item(21, 11).Value = 1.290846183014366
$ws.Cells.Item(21, 14).Value = 2.481744821705561
$ws.Cells.Item(22, 2).Value = 1.013570527292131
$ws.Cells.Item(22, 3).Value = 0.3402588151404871
$ws.Cells.Item(22, 4).Value = 0.09683822632146644
$ws.Cells.Item(22, 5).Value = 0.1550609635158722
$ws.Cells.Item(22, 6).Value = 2.834120401188784
$ws.Cells.Item(22, 9).Value = 1.387721642510073
$ws.Cells.Item(22, 10).Value = 0.2435712402594277
$ws.Cells.Item(22, 11).Value = 1.341201654321907
$ws.Cells.Item(22, 14).Value = 2.471178829844632
$ws.Cells.Item(23, 2).Value = 0.9922460700466331
$ws.Cells.Item(23, 3).Value = 0.3342475922470101
$ws.Cells.Item(23, 4).Value = 0.09546077353724058
$ws.Cells.Item(23, 5).Value = 0.1532785921244724
$ws.Cells.Item(23, 6).Value = 2.823666573541004
$ws.Cells.Item(23, 9).Value = 1.385438803817422
$ws.Cells.Item(23, 10).Value = 0.2410605890287485
$ws.Cells.Item(23, 11).Value = 1.314276732149068
$ws.Cells.Item(23, 14).Value = 2.476759230073654
$ws.Cells.Item(24, 2).Value = 0.9122625791175096
$ws.Cells.Item(24, 3).Value = 0.3117848052524153
$ws.Cells.Item(24, 4).Value = 0.09031497607581684
$ws.Cells.Item(24, 5).Value = 0.1466610663278942
$ws.Cells.Item(24, 6).Value = 2.786414162706507
$ws.Cells.Item(24, 9).Value = 1.37786952116474
$ws.Cells.Item(24, 10).Value = 0.2317738914086931
$ws.Cells.Item(24, 11).Value = 1.213384133365366
$ws.Cells.Item(24, 14).Value = 2.499185808514738
$ws.Cells.Item(25, 2).Value = 0.8276527040744952
$ws.Cells.Item(25, 3).Value = 0.2881991530478558
$ws.Cells.Item(25, 4).Value = 0.08491751243862211
$ws.Cells.Item(25, 5).Value = 0.1398065722140913
$ws.Cells.Item(25, 6).Value = 2.751214431737125
$ws.Cells.Item(25, 9).Value = 1.37199713150563
$ws.Cells.Item(25, 10).Value = 0.2222286259927557
$ws.Cells.Item(25, 11).Value = 1.106859545683477
$ws.Cells.Item(25, 14).Value = 2.526126182016057
